# Regenerate save_data to use K (strikeouts) instead of Strike# for column G,
# writing the recalculated s_vals into the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 3
    4  = 2
    5  = 5
    6  = 6
    7  = 2
    8  = 3
    9  = 6
    10 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
